$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

# Update MaxProd (F) and InvestCost (I) values for rows 8-18
$ws.Range("F8").Value = 18
$ws.Range("I8").Value = 8

$ws.Range("I9").Value = 8

$ws.Range("F10").Value = 16
$ws.Range("I10").Value = 8

$ws.Range("I11").Value = 8

$ws.Range("F12").Value = 15
$ws.Range("I12").Value = 8

$ws.Range("I13").Value = 8

$ws.Range("F14").Value = 6
$ws.Range("I14").Value = 8

$ws.Range("I15").Value = 8

$ws.Range("F16").Value = 80
$ws.Range("I16").Value = 8

$ws.Range("I17").Value = 8

$ws.Range("I18").Value = 8

# Update the active cell selection to I21
$ws.Activate()
$ws.Range("I21").Select()
